$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(61, 1).Value = "大智 (稳健智远)"
$ws.Cells.Item(61, 2).NumberFormat = "@"
$ws.Cells.Item(61, 2).Value = "000333"
$ws.Cells.Item(61, 2).ClearFormats()
$ws.Cells.Item(61, 3).Value = "美的集团"
$ws.Cells.Item(61, 4).Value = 3
$ws.Cells.Item(61, 5).Value = 42.89719527444501
$ws.Cells.Item(61, 6).Value = 71.26
$ws.Cells.Item(61, 7).Value = 3056.854135256952
$ws.Cells.Item(61, 8).Value = 101895.1378418984
$ws.Cells.Item(61, 9).NumberFormat = "@"
$ws.Cells.Item(61, 9).Value = "202506161355"
$ws.Cells.Item(61, 9).ClearFormats()
$ws.Cells.Item(62, 1).Value = "大智 (稳健智远)"
$ws.Cells.Item(62, 2).NumberFormat = "@"
$ws.Cells.Item(62, 2).Value = "510050"
$ws.Cells.Item(62, 2).ClearFormats()
$ws.Cells.Item(62, 3).Value = "上证50ETF"
$ws.Cells.Item(62, 4).Value = 5
$ws.Cells.Item(62, 5).Value = 1852.638869852698
$ws.Cells.Item(62, 6).Value = 2.75
$ws.Cells.Item(62, 7).Value = 5094.75689209492
$ws.Cells.Item(62, 8).Value = 101895.1378418984
$ws.Cells.Item(62, 9).NumberFormat = "@"
$ws.Cells.Item(62, 9).Value = "202506161355"
$ws.Cells.Item(62, 9).ClearFormats()
$ws.Cells.Item(63, 1).Value = "大智 (稳健智远)"
$ws.Cells.Item(63, 2).NumberFormat = "@"
$ws.Cells.Item(63, 2).Value = "510300"
$ws.Cells.Item(63, 2).ClearFormats()
$ws.Cells.Item(63, 3).Value = "沪深300ETF"
$ws.Cells.Item(63, 4).Value = 5
$ws.Cells.Item(63, 5).Value = 1276.881426590205
$ws.Cells.Item(63, 6).Value = 3.99
$ws.Cells.Item(63, 7).Value = 5094.75689209492
$ws.Cells.Item(63, 8).Value = 101895.1378418984
$ws.Cells.Item(63, 9).NumberFormat = "@"
$ws.Cells.Item(63, 9).Value = "202506161355"
$ws.Cells.Item(63, 9).ClearFormats()
$ws.Cells.Item(64, 1).Value = "大智 (稳健智远)"
$ws.Cells.Item(64, 2).NumberFormat = "@"
$ws.Cells.Item(64, 2).Value = "518880"
$ws.Cells.Item(64, 2).ClearFormats()
$ws.Cells.Item(64, 3).Value = "黄金ETF"
$ws.Cells.Item(64, 4).Value = 5
$ws.Cells.Item(64, 5).Value = 673.9096418114974
$ws.Cells.Item(64, 6).Value = 7.56
$ws.Cells.Item(64, 7).Value = 5094.75689209492
$ws.Cells.Item(64, 8).Value = 101895.1378418984
$ws.Cells.Item(64, 9).NumberFormat = "@"
$ws.Cells.Item(64, 9).Value = "202506161355"
$ws.Cells.Item(64, 9).ClearFormats()
$ws.Cells.Item(65, 1).Value = "大智 (稳健智远)"
$ws.Cells.Item(65, 2).NumberFormat = "@"
$ws.Cells.Item(65, 2).Value = "600085"
$ws.Cells.Item(65, 2).ClearFormats()
$ws.Cells.Item(65, 3).Value = "同仁堂"
$ws.Cells.Item(65, 4).Value = 2
$ws.Cells.Item(65, 5).Value = 56.06334956913254
$ws.Cells.Item(65, 6).Value = 36.35
$ws.Cells.Item(65, 7).Value = 2037.902756837968
$ws.Cells.Item(65, 8).Value = 101895.1378418984
$ws.Cells.Item(65, 9).NumberFormat = "@"
$ws.Cells.Item(65, 9).Value = "202506161355"
$ws.Cells.Item(65, 9).ClearFormats()
$ws.Cells.Item(66, 1).Value = "大智 (稳健智远)"
$ws.Cells.Item(66, 2).NumberFormat = "@"
$ws.Cells.Item(66, 2).Value = "600900"
$ws.Cells.Item(66, 2).ClearFormats()
$ws.Cells.Item(66, 3).Value = "长江电力"
$ws.Cells.Item(66, 4).Value = 20
$ws.Cells.Item(66, 5).Value = 669.0422707938175
$ws.Cells.Item(66, 6).Value = 30.46
$ws.Cells.Item(66, 7).Value = 20379.02756837968
$ws.Cells.Item(66, 8).Value = 101895.1378418984
$ws.Cells.Item(66, 9).NumberFormat = "@"
$ws.Cells.Item(66, 9).Value = "202506161355"
$ws.Cells.Item(66, 9).ClearFormats()
$ws.Cells.Item(67, 1).Value = "大智 (稳健智远)"
$ws.Cells.Item(67, 2).NumberFormat = "@"
$ws.Cells.Item(67, 2).Value = "600989"
$ws.Cells.Item(67, 2).ClearFormats()
$ws.Cells.Item(67, 3).Value = "宝丰能源"
$ws.Cells.Item(67, 4).Value = 5
$ws.Cells.Item(67, 5).Value = 308.5861230826723
$ws.Cells.Item(67, 6).Value = 16.51
$ws.Cells.Item(67, 7).Value = 5094.75689209492
$ws.Cells.Item(67, 8).Value = 101895.1378418984
$ws.Cells.Item(67, 9).NumberFormat = "@"
$ws.Cells.Item(67, 9).Value = "202506161355"
$ws.Cells.Item(67, 9).ClearFormats()
$ws.Cells.Item(68, 1).Value = "大智 (稳健智远)"
$ws.Cells.Item(68, 2).NumberFormat = "@"
$ws.Cells.Item(68, 2).Value = "601899"
$ws.Cells.Item(68, 2).ClearFormats()
$ws.Cells.Item(68, 3).Value = "XD紫金矿"
$ws.Cells.Item(68, 4).Value = 10
$ws.Cells.Item(68, 5).Value = 541.4194359293219
$ws.Cells.Item(68, 6).Value = 18.82
$ws.Cells.Item(68, 7).Value = 10189.51378418984
$ws.Cells.Item(68, 8).Value = 101895.1378418984
$ws.Cells.Item(68, 9).NumberFormat = "@"
$ws.Cells.Item(68, 9).Value = "202506161355"
$ws.Cells.Item(68, 9).ClearFormats()
$ws.Cells.Item(69, 1).Value = "大智 (稳健智远)"
$ws.Cells.Item(69, 2).NumberFormat = "@"
$ws.Cells.Item(69, 2).Value = "HK02899"
$ws.Cells.Item(69, 2).ClearFormats()
$ws.Cells.Item(69, 3).Value = "紫金矿业"
$ws.Cells.Item(69, 4).Value = 10
$ws.Cells.Item(69, 5).Value = 521.4694874201556
$ws.Cells.Item(69, 6).Value = 19.54
$ws.Cells.Item(69, 7).Value = 10189.51378418984
$ws.Cells.Item(69, 8).Value = 101895.1378418984
$ws.Cells.Item(69, 9).NumberFormat = "@"
$ws.Cells.Item(69, 9).Value = "202506161355"
$ws.Cells.Item(69, 9).ClearFormats()
$ws.Cells.Item(70, 1).Value = "大智 (稳健智远)"
$ws.Cells.Item(70, 2).NumberFormat = "@"
$ws.Cells.Item(70, 2).Value = "HK06881"
$ws.Cells.Item(70, 2).ClearFormats()
$ws.Cells.Item(70, 3).Value = "中国银河"
$ws.Cells.Item(70, 4).Value = 5
$ws.Cells.Item(70, 5).Value = 600.7968033130801
$ws.Cells.Item(70, 6).Value = 8.48
$ws.Cells.Item(70, 7).Value = 5094.75689209492
$ws.Cells.Item(70, 8).Value = 101895.1378418984
$ws.Cells.Item(70, 9).NumberFormat = "@"
$ws.Cells.Item(70, 9).Value = "202506161355"
$ws.Cells.Item(70, 9).ClearFormats()
$ws.Cells.Item(71, 1).Value = "大智 (稳健智远)"
$ws.Cells.Item(71, 2).NumberFormat = "@"
$ws.Cells.Item(71, 2).Value = "100000"
$ws.Cells.Item(71, 2).ClearFormats()
$ws.Cells.Item(71, 3).Value = "现金"
$ws.Cells.Item(71, 4).Value = 30
$ws.Cells.Item(71, 5).Value = 30568.54135256952
$ws.Cells.Item(71, 6).Value = 1
$ws.Cells.Item(71, 7).Value = 30568.54135256952
$ws.Cells.Item(71, 8).Value = 101895.1378418984
$ws.Cells.Item(71, 9).NumberFormat = "@"
$ws.Cells.Item(71, 9).Value = "202506161355"
$ws.Cells.Item(71, 9).ClearFormats()

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(44, 1).Value = "大成 (锐进先锋)"
$ws.Cells.Item(44, 2).NumberFormat = "@"
$ws.Cells.Item(44, 2).Value = "000725"
$ws.Cells.Item(44, 2).ClearFormats()
$ws.Cells.Item(44, 3).Value = "京东方A"
$ws.Cells.Item(44, 4).Value = 5
$ws.Cells.Item(44, 5).Value = 1264.047242468624
$ws.Cells.Item(44, 6).Value = 3.89
$ws.Cells.Item(44, 7).Value = 4917.143773202947
$ws.Cells.Item(44, 8).Value = 98342.87546405893
$ws.Cells.Item(44, 9).NumberFormat = "@"
$ws.Cells.Item(44, 9).Value = "202506161355"
$ws.Cells.Item(44, 9).ClearFormats()
$ws.Cells.Item(45, 1).Value = "大成 (锐进先锋)"
$ws.Cells.Item(45, 2).NumberFormat = "@"
$ws.Cells.Item(45, 2).Value = "159781"
$ws.Cells.Item(45, 2).ClearFormats()
$ws.Cells.Item(45, 3).Value = "科创创业ETF"
$ws.Cells.Item(45, 4).Value = 5
$ws.Cells.Item(45, 5).Value = 9277.629760760277
$ws.Cells.Item(45, 6).Value = 0.53
$ws.Cells.Item(45, 7).Value = 4917.143773202947
$ws.Cells.Item(45, 8).Value = 98342.87546405893
$ws.Cells.Item(45, 9).NumberFormat = "@"
$ws.Cells.Item(45, 9).Value = "202506161355"
$ws.Cells.Item(45, 9).ClearFormats()
$ws.Cells.Item(46, 1).Value = "大成 (锐进先锋)"
$ws.Cells.Item(46, 2).NumberFormat = "@"
$ws.Cells.Item(46, 2).Value = "513100"
$ws.Cells.Item(46, 2).ClearFormats()
$ws.Cells.Item(46, 3).Value = "纳指ETF"
$ws.Cells.Item(46, 4).Value = 5
$ws.Cells.Item(46, 5).Value = 3131.938709046463
$ws.Cells.Item(46, 6).Value = 1.57
$ws.Cells.Item(46, 7).Value = 4917.143773202947
$ws.Cells.Item(46, 8).Value = 98342.87546405893
$ws.Cells.Item(46, 9).NumberFormat = "@"
$ws.Cells.Item(46, 9).Value = "202506161355"
$ws.Cells.Item(46, 9).ClearFormats()
$ws.Cells.Item(47, 1).Value = "大成 (锐进先锋)"
$ws.Cells.Item(47, 2).NumberFormat = "@"
$ws.Cells.Item(47, 2).Value = "513290"
$ws.Cells.Item(47, 2).ClearFormats()
$ws.Cells.Item(47, 3).Value = "纳指生物科技ETF"
$ws.Cells.Item(47, 4).Value = 1
$ws.Cells.Item(47, 5).Value = 870.2909333102562
$ws.Cells.Item(47, 6).Value = 1.13
$ws.Cells.Item(47, 7).Value = 983.4287546405893
$ws.Cells.Item(47, 8).Value = 98342.87546405893
$ws.Cells.Item(47, 9).NumberFormat = "@"
$ws.Cells.Item(47, 9).Value = "202506161355"
$ws.Cells.Item(47, 9).ClearFormats()
$ws.Cells.Item(48, 1).Value = "大成 (锐进先锋)"
$ws.Cells.Item(48, 2).NumberFormat = "@"
$ws.Cells.Item(48, 2).Value = "603119"
$ws.Cells.Item(48, 2).ClearFormats()
$ws.Cells.Item(48, 3).Value = "浙江荣泰"
$ws.Cells.Item(48, 4).Value = 45
$ws.Cells.Item(48, 5).Value = 1069.978093782073
$ws.Cells.Item(48, 6).Value = 41.36
$ws.Cells.Item(48, 7).Value = 44254.29395882652
$ws.Cells.Item(48, 8).Value = 98342.87546405893
$ws.Cells.Item(48, 9).NumberFormat = "@"
$ws.Cells.Item(48, 9).Value = "202506161355"
$ws.Cells.Item(48, 9).ClearFormats()
$ws.Cells.Item(49, 1).Value = "大成 (锐进先锋)"
$ws.Cells.Item(49, 2).NumberFormat = "@"
$ws.Cells.Item(49, 2).Value = "688290"
$ws.Cells.Item(49, 2).ClearFormats()
$ws.Cells.Item(49, 3).Value = "景业智能"
$ws.Cells.Item(49, 4).Value = 9
$ws.Cells.Item(49, 5).Value = 147.121987895035
$ws.Cells.Item(49, 6).Value = 60.16
$ws.Cells.Item(49, 7).Value = 8850.858791765304
$ws.Cells.Item(49, 8).Value = 98342.87546405893
$ws.Cells.Item(49, 9).NumberFormat = "@"
$ws.Cells.Item(49, 9).Value = "202506161355"
$ws.Cells.Item(49, 9).ClearFormats()
$ws.Cells.Item(50, 1).Value = "大成 (锐进先锋)"
$ws.Cells.Item(50, 2).NumberFormat = "@"
$ws.Cells.Item(50, 2).Value = "100000"
$ws.Cells.Item(50, 2).ClearFormats()
$ws.Cells.Item(50, 3).Value = "现金"
$ws.Cells.Item(50, 4).Value = 30
$ws.Cells.Item(50, 5).Value = 29502.86263921768
$ws.Cells.Item(50, 6).Value = 1
$ws.Cells.Item(50, 7).Value = 29502.86263921768
$ws.Cells.Item(50, 8).Value = 98342.87546405893
$ws.Cells.Item(50, 9).NumberFormat = "@"
$ws.Cells.Item(50, 9).Value = "202506161355"
$ws.Cells.Item(50, 9).ClearFormats()

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(91, 1).Value = "范式进化投资组合"
$ws.Cells.Item(91, 2).NumberFormat = "@"
$ws.Cells.Item(91, 2).Value = "000333"
$ws.Cells.Item(91, 2).ClearFormats()
$ws.Cells.Item(91, 3).Value = "美的集团"
$ws.Cells.Item(91, 4).Value = 1
$ws.Cells.Item(91, 5).Value = 14.02515882310653
$ws.Cells.Item(91, 6).Value = 71.26
$ws.Cells.Item(91, 7).Value = 999.4328177345716
$ws.Cells.Item(91, 8).Value = 99943.28177345716
$ws.Cells.Item(91, 9).NumberFormat = "@"
$ws.Cells.Item(91, 9).Value = "202506161355"
$ws.Cells.Item(91, 9).ClearFormats()
$ws.Cells.Item(92, 1).Value = "范式进化投资组合"
$ws.Cells.Item(92, 2).NumberFormat = "@"
$ws.Cells.Item(92, 2).Value = "000725"
$ws.Cells.Item(92, 2).ClearFormats()
$ws.Cells.Item(92, 3).Value = "京东方A"
$ws.Cells.Item(92, 4).Value = 5
$ws.Cells.Item(92, 5).Value = 1284.618017653691
$ws.Cells.Item(92, 6).Value = 3.89
$ws.Cells.Item(92, 7).Value = 4997.164088672858
$ws.Cells.Item(92, 8).Value = 99943.28177345716
$ws.Cells.Item(92, 9).NumberFormat = "@"
$ws.Cells.Item(92, 9).Value = "202506161355"
$ws.Cells.Item(92, 9).ClearFormats()
$ws.Cells.Item(93, 1).Value = "范式进化投资组合"
$ws.Cells.Item(93, 2).NumberFormat = "@"
$ws.Cells.Item(93, 2).Value = "159781"
$ws.Cells.Item(93, 2).ClearFormats()
$ws.Cells.Item(93, 3).Value = "科创创业ETF"
$ws.Cells.Item(93, 4).Value = 5
$ws.Cells.Item(93, 5).Value = 9428.611488061997
$ws.Cells.Item(93, 6).Value = 0.53
$ws.Cells.Item(93, 7).Value = 4997.164088672858
$ws.Cells.Item(93, 8).Value = 99943.28177345716
$ws.Cells.Item(93, 9).NumberFormat = "@"
$ws.Cells.Item(93, 9).Value = "202506161355"
$ws.Cells.Item(93, 9).ClearFormats()
$ws.Cells.Item(94, 1).Value = "范式进化投资组合"
$ws.Cells.Item(94, 2).NumberFormat = "@"
$ws.Cells.Item(94, 2).Value = "510050"
$ws.Cells.Item(94, 2).ClearFormats()
$ws.Cells.Item(94, 3).Value = "上证50ETF"
$ws.Cells.Item(94, 4).Value = 5
$ws.Cells.Item(94, 5).Value = 1817.150577699221
$ws.Cells.Item(94, 6).Value = 2.75
$ws.Cells.Item(94, 7).Value = 4997.164088672858
$ws.Cells.Item(94, 8).Value = 99943.28177345716
$ws.Cells.Item(94, 9).NumberFormat = "@"
$ws.Cells.Item(94, 9).Value = "202506161355"
$ws.Cells.Item(94, 9).ClearFormats()
$ws.Cells.Item(95, 1).Value = "范式进化投资组合"
$ws.Cells.Item(95, 2).NumberFormat = "@"
$ws.Cells.Item(95, 2).Value = "510300"
$ws.Cells.Item(95, 2).ClearFormats()
$ws.Cells.Item(95, 3).Value = "沪深300ETF"
$ws.Cells.Item(95, 4).Value = 5
$ws.Cells.Item(95, 5).Value = 1252.422077361618
$ws.Cells.Item(95, 6).Value = 3.99
$ws.Cells.Item(95, 7).Value = 4997.164088672858
$ws.Cells.Item(95, 8).Value = 99943.28177345716
$ws.Cells.Item(95, 9).NumberFormat = "@"
$ws.Cells.Item(95, 9).Value = "202506161355"
$ws.Cells.Item(95, 9).ClearFormats()
$ws.Cells.Item(96, 1).Value = "范式进化投资组合"
$ws.Cells.Item(96, 2).NumberFormat = "@"
$ws.Cells.Item(96, 2).Value = "513100"
$ws.Cells.Item(96, 2).ClearFormats()
$ws.Cells.Item(96, 3).Value = "纳指ETF"
$ws.Cells.Item(96, 4).Value = 1
$ws.Cells.Item(96, 5).Value = 636.5814125697908
$ws.Cells.Item(96, 6).Value = 1.57
$ws.Cells.Item(96, 7).Value = 999.4328177345716
$ws.Cells.Item(96, 8).Value = 99943.28177345716
$ws.Cells.Item(96, 9).NumberFormat = "@"
$ws.Cells.Item(96, 9).Value = "202506161355"
$ws.Cells.Item(96, 9).ClearFormats()
$ws.Cells.Item(97, 1).Value = "范式进化投资组合"
$ws.Cells.Item(97, 2).NumberFormat = "@"
$ws.Cells.Item(97, 2).Value = "513290"
$ws.Cells.Item(97, 2).ClearFormats()
$ws.Cells.Item(97, 3).Value = "纳指生物科技ETF"
$ws.Cells.Item(97, 4).Value = 1
$ws.Cells.Item(97, 5).Value = 884.4538210040457
$ws.Cells.Item(97, 6).Value = 1.13
$ws.Cells.Item(97, 7).Value = 999.4328177345716
$ws.Cells.Item(97, 8).Value = 99943.28177345716
$ws.Cells.Item(97, 9).NumberFormat = "@"
$ws.Cells.Item(97, 9).Value = "202506161355"
$ws.Cells.Item(97, 9).ClearFormats()
$ws.Cells.Item(98, 1).Value = "范式进化投资组合"
$ws.Cells.Item(98, 2).NumberFormat = "@"
$ws.Cells.Item(98, 2).Value = "518880"
$ws.Cells.Item(98, 2).ClearFormats()
$ws.Cells.Item(98, 3).Value = "黄金ETF"
$ws.Cells.Item(98, 4).Value = 1
$ws.Cells.Item(98, 5).Value = 132.2001081659486
$ws.Cells.Item(98, 6).Value = 7.56
$ws.Cells.Item(98, 7).Value = 999.4328177345715
$ws.Cells.Item(98, 8).Value = 99943.28177345716
$ws.Cells.Item(98, 9).NumberFormat = "@"
$ws.Cells.Item(98, 9).Value = "202506161355"
$ws.Cells.Item(98, 9).ClearFormats()
$ws.Cells.Item(99, 1).Value = "范式进化投资组合"
$ws.Cells.Item(99, 2).NumberFormat = "@"
$ws.Cells.Item(99, 2).Value = "600085"
$ws.Cells.Item(99, 2).ClearFormats()
$ws.Cells.Item(99, 3).Value = "同仁堂"
$ws.Cells.Item(99, 4).Value = 1
$ws.Cells.Item(99, 5).Value = 27.49471300507762
$ws.Cells.Item(99, 6).Value = 36.35
$ws.Cells.Item(99, 7).Value = 999.4328177345716
$ws.Cells.Item(99, 8).Value = 99943.28177345716
$ws.Cells.Item(99, 9).NumberFormat = "@"
$ws.Cells.Item(99, 9).Value = "202506161355"
$ws.Cells.Item(99, 9).ClearFormats()
$ws.Cells.Item(100, 1).Value = "范式进化投资组合"
$ws.Cells.Item(100, 2).NumberFormat = "@"
$ws.Cells.Item(100, 2).Value = "600900"
$ws.Cells.Item(100, 2).ClearFormats()
$ws.Cells.Item(100, 3).Value = "长江电力"
$ws.Cells.Item(100, 4).Value = 1
$ws.Cells.Item(100, 5).Value = 32.81132034584936
$ws.Cells.Item(100, 6).Value = 30.46
$ws.Cells.Item(100, 7).Value = 999.4328177345716
$ws.Cells.Item(100, 8).Value = 99943.28177345716
$ws.Cells.Item(100, 9).NumberFormat = "@"
$ws.Cells.Item(100, 9).Value = "202506161355"
$ws.Cells.Item(100, 9).ClearFormats()
$ws.Cells.Item(101, 1).Value = "范式进化投资组合"
$ws.Cells.Item(101, 2).NumberFormat = "@"
$ws.Cells.Item(101, 2).Value = "600989"
$ws.Cells.Item(101, 2).ClearFormats()
$ws.Cells.Item(101, 3).Value = "宝丰能源"
$ws.Cells.Item(101, 4).Value = 5
$ws.Cells.Item(101, 5).Value = 302.6749902285196
$ws.Cells.Item(101, 6).Value = 16.51
$ws.Cells.Item(101, 7).Value = 4997.164088672858
$ws.Cells.Item(101, 8).Value = 99943.28177345716
$ws.Cells.Item(101, 9).NumberFormat = "@"
$ws.Cells.Item(101, 9).Value = "202506161355"
$ws.Cells.Item(101, 9).ClearFormats()
$ws.Cells.Item(102, 1).Value = "范式进化投资组合"
$ws.Cells.Item(102, 2).NumberFormat = "@"
$ws.Cells.Item(102, 2).Value = "601899"
$ws.Cells.Item(102, 2).ClearFormats()
$ws.Cells.Item(102, 3).Value = "XD紫金矿"
$ws.Cells.Item(102, 4).Value = 10
$ws.Cells.Item(102, 5).Value = 531.0482559694855
$ws.Cells.Item(102, 6).Value = 18.82
$ws.Cells.Item(102, 7).Value = 9994.328177345717
$ws.Cells.Item(102, 8).Value = 99943.28177345716
$ws.Cells.Item(102, 9).NumberFormat = "@"
$ws.Cells.Item(102, 9).Value = "202506161355"
$ws.Cells.Item(102, 9).ClearFormats()
$ws.Cells.Item(103, 1).Value = "范式进化投资组合"
$ws.Cells.Item(103, 2).NumberFormat = "@"
$ws.Cells.Item(103, 2).Value = "603119"
$ws.Cells.Item(103, 2).ClearFormats()
$ws.Cells.Item(103, 3).Value = "浙江荣泰"
$ws.Cells.Item(103, 4).Value = 1
$ws.Cells.Item(103, 5).Value = 24.16423640557475
$ws.Cells.Item(103, 6).Value = 41.36
$ws.Cells.Item(103, 7).Value = 999.4328177345716
$ws.Cells.Item(103, 8).Value = 99943.28177345716
$ws.Cells.Item(103, 9).NumberFormat = "@"
$ws.Cells.Item(103, 9).Value = "202506161355"
$ws.Cells.Item(103, 9).ClearFormats()
$ws.Cells.Item(104, 1).Value = "范式进化投资组合"
$ws.Cells.Item(104, 2).NumberFormat = "@"
$ws.Cells.Item(104, 2).Value = "HK06881"
$ws.Cells.Item(104, 2).ClearFormats()
$ws.Cells.Item(104, 3).Value = "中国银河"
$ws.Cells.Item(104, 4).Value = 1
$ws.Cells.Item(104, 5).Value = 117.8576436007749
$ws.Cells.Item(104, 6).Value = 8.48
$ws.Cells.Item(104, 7).Value = 999.4328177345716
$ws.Cells.Item(104, 8).Value = 99943.28177345716
$ws.Cells.Item(104, 9).NumberFormat = "@"
$ws.Cells.Item(104, 9).Value = "202506161355"
$ws.Cells.Item(104, 9).ClearFormats()
$ws.Cells.Item(105, 1).Value = "范式进化投资组合"
$ws.Cells.Item(105, 2).NumberFormat = "@"
$ws.Cells.Item(105, 2).Value = "100000"
$ws.Cells.Item(105, 2).ClearFormats()
$ws.Cells.Item(105, 3).Value = "现金"
$ws.Cells.Item(105, 4).Value = 57
$ws.Cells.Item(105, 5).Value = 56967.67061087058
$ws.Cells.Item(105, 6).Value = 1
$ws.Cells.Item(105, 7).Value = 56967.67061087058
$ws.Cells.Item(105, 8).Value = 99943.28177345716
$ws.Cells.Item(105, 9).NumberFormat = "@"
$ws.Cells.Item(105, 9).Value = "202506161355"
$ws.Cells.Item(105, 9).ClearFormats()
